$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.636.20"
$ws.Range("E2").Value = "  -0.57%  "
$ws.Range("D3").Value = "3.916.30"
$ws.Range("E3").Value = "  +3.28%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "602.98"
$ws.Range("E5").Value = "  +0.34%  "
$ws.Range("D6").Value = "165.74"
$ws.Range("E6").Value = "  +1.97%  "
$ws.Range("D7").Value = "3.914.78"
$ws.Range("E7").Value = "  +3.33%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("E9").Value = "  -1.08%  "
$ws.Range("E10").Value = "  -1.70%  "
$ws.Range("D11").Value = "6.38"
$ws.Range("E11").Value = "  +1.21%  "
$ws.Range("D12").Value = "'0.460"
$ws.Range("E12").Value = "  +0.35%  "
$ws.Range("D13").Value = "0.0000255"
$ws.Range("E13").Value = "  +4.31%  "
$ws.Range("D14").Value = "'37.30"
$ws.Range("E14").Value = "  +0.44%  "
$ws.Range("D15").Value = "4.572.62"
$ws.Range("E15").Value = "  +3.20%  "
$ws.Range("D16").Value = "3.925.46"
$ws.Range("E16").Value = "  +3.51%  "
$ws.Range("D17").Value = "68.742.20"
$ws.Range("E17").Value = "  -0.57%  "
$ws.Range("D18").Value = "7.45"
$ws.Range("E18").Value = "  +1.01%  "
$ws.Range("D19").Value = "17.13"
$ws.Range("E19").Value = "  -0.89%  "
$ws.Range("E20").Value = "  -2.14%  "
$ws.Range("E21").Value = "  -2.28%  "
$ws.Range("D22").Value = "485.81"
$ws.Range("D23").Value = "0.723"
$ws.Range("E23").Value = "  +0.55%  "
$ws.Range("D24").Value = "'0.0000170"
$ws.Range("E24").Value = "  +11.94%  "
$ws.Range("D25").Value = "84.47"
$ws.Range("E25").Value = "  -0.04%  "
$ws.Range("E26").Value = "  -0.39%  "
$ws.Range("D27").Value = "12.07"
$ws.Range("E27").Value = "  -0.86%  "
$ws.Range("D28").Value = "10.12"
$ws.Range("E28").Value = "  +0.79%  "
$ws.Range("E29").Value = "  +0.12%  "
$ws.Range("E30").Value = "  -1.03%  "
$ws.Range("D31").Value = "4.069.64"
$ws.Range("E31").Value = "  +3.19%  "
$ws.Range("D32").Value = "2.39"
$ws.Range("E32").Value = "  -0.11%  "
$ws.Range("D33").Value = "7.77"
$ws.Range("E33").Value = "  -3.27%  "
$ws.Range("D34").Value = "31.91"
$ws.Range("E34").Value = "  +0.25%  "
$ws.Range("D35").Value = "3.868.78"
$ws.Range("E35").Value = "  +3.39%  "
$ws.Range("E36").Value = "  +0.33%  "
$ws.Range("E37").Value = "  +2.66%  "
$ws.Range("D38").Value = "5.91"
$ws.Range("E38").Value = "  +0.59%  "
$ws.Range("E39").Value = "  -1.51%  "
$ws.Range("D40").Value = "3.18"
$ws.Range("E40").Value = "  +5.66%  "
$ws.Range("D41").Value = "0.999"
$ws.Range("E41").Value = "  -0.06%  "
$ws.Range("D42").Value = "0.314"
$ws.Range("E42").Value = "  -2.13%  "
$ws.Range("D43").Value = "430.29"
$ws.Range("E43").Value = "  +3.16%  "
$ws.Range("D44").Value = "48.45"
$ws.Range("E44").Value = "  -0.05%  "
$ws.Range("E45").Value = "  +0.01%  "
$ws.Range("D46").Value = "8.49"
$ws.Range("E46").Value = "  +1.56%  "
$ws.Range("E47").Value = "  +0.01%  "
$ws.Range("D48").Value = "26.32"
$ws.Range("E48").Value = "  +7.50%  "
$ws.Range("D49").Value = "141.72"
$ws.Range("E49").Value = "  +0.05%  "
$ws.Range("D50").Value = "2.813.24"
$ws.Range("E50").Value = "  -0.12%  "
$ws.Range("E51").Value = "  +0.66%  "
